$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("summary")
$ws2 = $wb.Worksheets.Item("model_fit")

# Sheet "summary": update per-item statistics (rows 2-8)
$ws1.Range("C2").Value = 706
$ws1.Range("D2").Value = 652
$ws1.Range("E2").Value = 85.58
$ws1.Range("F2").Value = -2.19
$ws1.Range("H2").Value = 1
$ws1.Range("I2").Value = 0.08
$ws1.Range("J2").Value = 0.25
$ws1.Range("L2").Value = 1.13

$ws1.Range("C3").Value = 706
$ws1.Range("D3").Value = 598
$ws1.Range("E3").Value = 58.7
$ws1.Range("F3").Value = -0.42
$ws1.Range("I3").Value = 0.44
$ws1.Range("J3").Value = 0.31
$ws1.Range("K3").Value = 0.05
$ws1.Range("L3").Value = 0.97

$ws1.Range("C4").Value = 706
$ws1.Range("D4").Value = 361
$ws1.Range("E4").Value = 24.93
$ws1.Range("F4").Value = 1.42
$ws1.Range("G4").Value = 0.13
$ws1.Range("H4").Value = 1.07
$ws1.Range("I4").Value = 1.08
$ws1.Range("L4").Value = 0.71

$ws1.Range("C5").Value = 706
$ws1.Range("D5").Value = 671
$ws1.Range("E5").Value = 39.64
$ws1.Range("F5").Value = 0.54
$ws1.Range("I5").Value = 0.36
$ws1.Range("J5").Value = 0.32
$ws1.Range("K5").Value = 0.08
$ws1.Range("L5").Value = 1.1

$ws1.Range("C6").Value = 706
$ws1.Range("D6").Value = 678
$ws1.Range("E6").Value = 37.46
$ws1.Range("F6").Value = 0.66
$ws1.Range("H6").Value = 1.02
$ws1.Range("I6").Value = 0.43
$ws1.Range("K6").Value = 0.07
$ws1.Range("L6").Value = 1.09

$ws1.Range("C7").Value = 706
$ws1.Range("D7").Value = 676
$ws1.Range("E7").Value = 26.48
$ws1.Range("F7").Value = 1.29
$ws1.Range("J7").Value = 0.38
$ws1.Range("K7").Value = 0.05
$ws1.Range("L7").Value = 1.81

$ws1.Range("C8").Value = 706
$ws1.Range("D8").Value = 676
$ws1.Range("E8").Value = 23.08
$ws1.Range("F8").Value = 1.53
$ws1.Range("I8").Value = -0.51
$ws1.Range("J8").Value = 0.33
$ws1.Range("K8").Value = 0.05
$ws1.Range("L8").Value = 1.47

# Sheet "model_fit": update model fit statistics (rows 2-3)
$ws2.Range("B2").Value = 706
$ws2.Range("D2").Value = 4811
$ws2.Range("E2").Value = 4827
$ws2.Range("F2").Value = 4864
$ws2.Range("G2").Value = 0.564
$ws2.Range("H2").Value = 0.399

$ws2.Range("B3").Value = 706
$ws2.Range("D3").Value = 4795
$ws2.Range("E3").Value = 4823
$ws2.Range("F3").Value = 4887
$ws2.Range("G3").Value = 0.578
$ws2.Range("H3").Value = 0.404
